# Generate Report for Handoff
# Updates "Latest Handoff Datetime" values for the two rows whose handoff
# was just (re)generated (2792beb7... and 43044a75...) on both the zh-cn
# and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D12").Value = "2016-03-08 18:27:36"
$wsZhCn.Range("D14").Value = "2016-03-08 18:27:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D12").Value = "2016-03-08 18:27:44"
$wsDeDe.Range("D14").Value = "2016-03-08 18:27:44"
